$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 (Instance) value from "test2" to "wild1"
$ws.Range("D2").Value = "wild1"

# Update the active selection to D2 (previously B3)
$ws.Range("D2").Select()
